# Insert a new price-list row for "Katilda, Dawnhart Martyr: Katilda's Rising
# Dawn" (Innistrad: Crimson Vow) just above the existing "Sunfall" row, which
# pushes every following row down by one and widens the trailing SUM formula
# to include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift row 23 (and everything below it) down one row, inheriting formatting
# from the row being displaced - this also keeps the final SUM(D2:D100)
# formula auto-expanding to SUM(D2:D101) since Excel treats it as a normal
# row insertion.
$ws.Rows("23:23").Insert()

$ws.Range("A23").Value = "Katilda, Dawnhart Martyr: Katilda's Rising Dawn"
$ws.Range("B23").Value = "Innistrad: Crimson Vow"
$ws.Range("C23").Value = "Normal"
$ws.Range("D23").Value = 2.32

# Mirror the author's final cursor position on the new row.
$ws.Range("D24").Select()
